# [ADD] Added pipeline to test symbolic dynamic models
#
# Fills in the previously-empty experiment rows (14-17) on Sheet1 with the
# new PD/PID controller test-pipeline entries, and updates the active
# selection to match the author's saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - PD controller test step response
$ws.Range("C14").Value = "PD controller test step response"
$ws.Range("D14").Value = "Control the system via the PD controller and set the theta`nreference to a constant value"
$ws.Range("E14").Value = "Check static performance of the`ncontroller"
$ws.Range("G14").Value = "Paolo"

# Row 15 - PD controller test dynamic response
$ws.Range("C15").Value = "PD controller test dynamic`nresponse"
$ws.Range("D15").Value = "Control the system via the PD controller and set the theta`nreference to a sinusoid slower than the cutting frequency"
$ws.Range("E15").Value = "Check dynamic performance of the`ncontroller"
$ws.Range("G15").Value = "Paolo"

# Row 16 - PID controller test step response
$ws.Range("C16").Value = "PID controller test step response"
$ws.Range("D16").Value = "Control the system via the PID controller and set the theta`nreference to a constant value"
$ws.Range("E16").Value = "Check static performance of the`ncontroller"
$ws.Range("G16").Value = "Paolo"

# Row 17 - PID controller test dynamic response
$ws.Range("C17").Value = "PID controller test dynamic`nresponse"
$ws.Range("D17").Value = "Control the system via the PID controller and set the theta`nreference to a sinusoid slower than the cutting frequency"
$ws.Range("E17").Value = "Check dynamic performance of the`ncontroller"
$ws.Range("G17").Value = "Paolo"

# Match the author's saved view state: scrolled down with D17 selected.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select() | Out-Null
